# Update FuelPrices at 2025-04-14 02:41
# - Change C23's number format from date-only (YYYY-MM-DD) to
#   date-time (YYYY-MM-DD HH:MM:SS), matching the rest of the column.
# - Append a new data row 24 with values 779.953 / 690.04 / 45755,
#   with C24 using the date-only (YYYY-MM-DD) format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the number format on the previous last row's date cell.
$ws.Range("C23").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new row of data.
$ws.Range("A24").Value = 779.953
$ws.Range("B24").Value = 690.04
$ws.Range("C24").Value = 45755
$ws.Range("C24").NumberFormat = "YYYY-MM-DD"
